$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "blue" review-confirmation value in row 6 from "confirm" to "no"
$ws.Range("G6").Value = "no"

# Move selection to reflect the authored state (cursor left on G7 after the edit)
$ws.Range("G7").Select() | Out-Null
